$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "directory" sheet: sort the B:E block (file name / date / PC-type /
#    number) by date ascending, leaving the A numbering column (1..54)
#    untouched, then tack on a new "col_names" column F.
# ---------------------------------------------------------------------
$dir = $wb.Worksheets.Item("directory")
$dir.Range("B2:E55").Sort($dir.Range("C2:C55"))

$dir.Range("F1").Value = "col_names"
$dir.Columns.Item(6).ColumnWidth = 9.8

$dir.Range("C2").Select()

# ---------------------------------------------------------------------
# 2. New sheet "column_name_rename": a small lookup block in column B.
# ---------------------------------------------------------------------
$rename1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$rename1.Name = "column_name_rename"

$rename1.Range("B1").Value = 40513
$rename1.Range("B1").NumberFormat = "mmm-yy"

$rename1Values = @(
    "metro/rest_of_state",
    "postcode",
    "1BR_flats_counts",
    "1BR_flats_median",
    "2BR_flats_counts",
    "2BR_flats_median"
)
for ($i = 0; $i -lt $rename1Values.Count; $i++) {
    $rename1.Cells.Item($i + 2, 2).Value = $rename1Values[$i]
}

$rename1.Range("B5").Select()

# ---------------------------------------------------------------------
# 3. New sheet "rename_2011": old_name -> new_name mapping table.
# ---------------------------------------------------------------------
$rename2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$rename2.Name = "rename_2011"

$mapping = @(
    @("old_name", "new_name"),
    @("Metro/Rest of State", "metro_rest_of_state"),
    @("Postcode", "postcode"),
    @("1 BR  Flats", "1BR_flats_counts"),
    @("...4", "1BR_flats_median"),
    @("2 BR Flats", "2BR_flats_counts"),
    @("...6", "2BR_flats_median"),
    @("3 BR Flats", "3BR_flats_counts"),
    @("...8", "3BR_flats_median"),
    @("4+ BR Flats", "4BR_flats_counts"),
    @("...10", "4BR_flats_median"),
    @("Flat unknown BR", "flat_unknown_BR_counts"),
    @("...12", "flat_unknow_BR_median"),
    @("Total Flats", "total_flats_counts"),
    @("...14", "total_flats_median"),
    @("1 BR Houses", "1BR_houses_counts"),
    @("...16", "1BR_houses_median"),
    @("2 BR Houses", "2BR_houses_counts"),
    @("...18", "2BR_houses_median"),
    @("3 BR Houses", "3BR_houses_counts"),
    @("...20", "3BR_houses_median"),
    @("4+ BR Houses", "4BR_houses_counts"),
    @("...22", "4BR_houses_median"),
    @("House unknown BR", "house_unknown_BR_counts"),
    @("...24", "house_unknow_BR_median"),
    @("Total Houses", "total_houses_counts"),
    @("...26", "total_houses_median"),
    @("Other/Unknown", "Other_unknown_counts"),
    @("...28", "Other_unknown_median"),
    @("Total", "Total_count"),
    @("...30", "Total_unknown"),
    @("Quarter", "quarter")
)
for ($i = 0; $i -lt $mapping.Count; $i++) {
    $rename2.Cells.Item($i + 1, 1).Value = $mapping[$i][0]
    $rename2.Cells.Item($i + 1, 2).Value = $mapping[$i][1]
}

$rename2.Columns.Item(1).ColumnWidth = 23.36328125

$rename2.Activate()
$rename2.Range("H32:J33").Select()
